$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '61.633.93'
$ws.Range("E2").Value = '  -3.64%  '

# Row 3
$ws.Range("D3").Value = '2.478.19'
$ws.Range("E3").Value = '  -6.39%  '

# Row 4
$ws.Range("E4").Value = '  +0.08%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '555.51'
$ws.Range("E5").Value = '  -4.44%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '147.82'
$ws.Range("E6").Value = '  -5.55%  '

# Row 7
$ws.Range("E7").Value = '  +0.06%  '

# Row 8
$ws.Range("E8").Value = '  -3.43%  '

# Row 9
$ws.Range("D9").Value = '2.476.15'
$ws.Range("E9").Value = '  -6.47%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.108'
$ws.Range("E10").Value = '  -8.54%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.50'
$ws.Range("E11").Value = '  -5.50%  '

# Row 12
$ws.Range("E12").Value = '  -1.43%  '

# Row 13
$ws.Range("E13").Value = '  -6.61%  '

# Row 14
$ws.Range("E14").Value = '  -7.38%  '

# Row 15
$ws.Range("D15").Value = '2.928.78'
$ws.Range("E15").Value = '  -6.20%  '

# Row 16
$ws.Range("E16").Value = '  -8.68%  '

# Row 17
$ws.Range("D17").Value = '61.542.36'
$ws.Range("E17").Value = '  -3.65%  '

# Row 18
$ws.Range("D18").Value = '2.481.19'
$ws.Range("E18").Value = '  -6.47%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.25'
$ws.Range("E19").Value = '  -7.82%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.16'
$ws.Range("E20").Value = '  -7.43%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.24'
$ws.Range("E21").Value = '  -6.86%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '322.54'
$ws.Range("E22").Value = '  -6.91%  '

# Row 23
$ws.Range("E23").Value = '  -0.02%  '

# Row 24
$ws.Range("E24").Value = '  +2.20%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '64.45'
$ws.Range("E25").Value = '  -5.45%  '

# Row 26
$ws.Range("E26").Value = '  -9.54%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '563.89'
$ws.Range("E27").Value = '  -3.92%  '

# Row 28
$ws.Range("D28").Value = '2.606.63'
$ws.Range("E28").Value = '  -6.21%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.51'
$ws.Range("E29").Value = '  -6.39%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.999'
$ws.Range("E30").Value = '  +0.00%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.86'
$ws.Range("E31").Value = '  -4.86%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.34'
$ws.Range("E32").Value = '  -10.40%  '

# Row 33
$ws.Range("E33").Value = '  -6.66%  '

# Row 34
$ws.Range("E34").Value = '  -6.15%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.59'
$ws.Range("E35").Value = '  -8.58%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.93'
$ws.Range("E36").Value = '  -10.60%  '

# Row 37
$ws.Range("E37").Value = '  -10.76%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.00'
$ws.Range("E38").Value = '  +0.12%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.385'
$ws.Range("E39").Value = '  -4.95%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.62'
$ws.Range("E40").Value = '  -5.80%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '146.35'
$ws.Range("E41").Value = '  -3.25%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.76'
$ws.Range("E42").Value = '  -8.70%  '

# Row 43
$ws.Range("E43").Value = '  +0.13%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '40.61'
$ws.Range("E44").Value = '  -3.12%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.43'
$ws.Range("E45").Value = '  -5.74%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '148.97'
$ws.Range("E46").Value = '  -8.66%  '

# Row 47
$ws.Range("B47").Value = 'Filecoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.66'
$ws.Range("E47").Value = '  -6.65%  '

# Row 48
$ws.Range("B48").Value = 'InjectiveProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '22.06'
$ws.Range("E48").Value = '  -9.81%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0543'
$ws.Range("E49").Value = '  -8.10%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.597'
$ws.Range("E50").Value = '  -6.16%  '

# Row 51
$ws.Range("E51").Value = '  -5.45%  '
